# Fruta / hortaliza, semanal
#
# A new weekly price record (2022-06-03, serial 44715) is inserted as the
# new row 10 of the data table. All the rows that used to be 10-36 shift
# down by one (to 11-37), and the sheet's used range grows from
# A1:R36 to A1:R37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 10 (and everything below it) down by one row.
$ws.Rows("10").Insert()

# Fill in the brand-new row 10 with this week's record.
$ws.Cells.Item(10, 1).Value  = 3
$ws.Cells.Item(10, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(10, 3).Value  = "Coquimbo"
$ws.Cells.Item(10, 4).Value  = 44715
$ws.Cells.Item(10, 5).Value  = 5
$ws.Cells.Item(10, 6).Value  = 100112035
$ws.Cells.Item(10, 7).Value  = "Bruselas (repollito)"
$ws.Cells.Item(10, 8).Value  = "Sin especificar"
$ws.Cells.Item(10, 9).Value  = "Primera"
$ws.Cells.Item(10, 10).Value = 85
$ws.Cells.Item(10, 11).Value = 15000
$ws.Cells.Item(10, 12).Value = 15500
$ws.Cells.Item(10, 13).Value = 15235
$ws.Cells.Item(10, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(10, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(10, 16).Value = 1016
$ws.Cells.Item(10, 17).Value = 15
$ws.Cells.Item(10, 18).Value = "Hortaliza"
